$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Model")

# Row 2 (person_id): drop the "from CDM PERSONS" note.
$ws.Range("E2").ClearContents()

# Row 3 (birth_date): drop description + note, keep Varname/Format/Calculated.
$ws.Range("B3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Rows.Item(3).AutoFit()

# Row 4 (death_date): drop description + note, keep Varname/Format/Calculated.
$ws.Range("B4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Rows.Item(4).AutoFit()

# Row 5: replace sex_at_instance_creation with the new gender variable.
$ws.Range("A5").Value = "gender"
$ws.Range("B5").Clear()
$ws.Range("D5").Value = "M" + [char]10 + "F"
$ws.Range("D5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 31.2
$ws.Range("E5").ClearContents()

# Rows 6-9 (birth_month_imputed, birth_day_imputed, death_month_imputed,
# death_day_imputed): these parameters are removed entirely.
$ws.Range("A6").Clear()
$ws.Range("C6").Clear()
$ws.Range("I6").Clear()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()

$ws.Range("A7").Clear()
$ws.Range("C7").Clear()
$ws.Range("I7").Clear()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Rows.Item(7).AutoFit()

$ws.Range("A8").Clear()
$ws.Range("C8").Clear()
$ws.Range("I8").Clear()
$ws.Range("D8").ClearContents()
$ws.Rows.Item(8).AutoFit()

$ws.Range("A9").Clear()
$ws.Range("C9").Clear()
$ws.Range("I9").Clear()
$ws.Range("D9").ClearContents()
$ws.Rows.Item(9).AutoFit()
